$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-145 down to 67-146.
$ws.Rows("66:66").Insert()

# Populate the newly inserted row 66 with the new record's data.
$ws.Cells.Item(66, 1).Value = 10
$ws.Cells.Item(66, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(66, 3).Value = "La Araucanía"
$ws.Cells.Item(66, 4).Value = 45174
$ws.Cells.Item(66, 5).Value = 9
$ws.Cells.Item(66, 6).Value = 100112010
$ws.Cells.Item(66, 7).Value = "Achicoria"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 100
$ws.Cells.Item(66, 11).Value = 10000
$ws.Cells.Item(66, 12).Value = 10000
$ws.Cells.Item(66, 13).Value = 10000
$ws.Cells.Item(66, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(66, 15).Value = "Región Metropolitana"
$ws.Cells.Item(66, 16).Value = 556
$ws.Cells.Item(66, 17).Value = 18
$ws.Cells.Item(66, 18).Value = "Hortaliza"
